$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text / URL cells (Coin name, Link, Volume label columns) - plain text, no special handling needed
$textUpdates = @{
    "B10" = "One"
    "C10" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E10" = "9OneONEBestin24h"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E11" = "10WazirXWRX"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E12" = "11MandalaExchangeTokenMDX"
    "B13" = "LiechtensteinCryptoassetsExchange"
    "C13" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E13" = "12LiechtensteinCryptoassetsExchangeLCX"
    "B19" = "TigerCash"
    "C19" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "E19" = "18TigerCashTCH"
    "B20" = "HotbitToken"
    "C20" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "E20" = "19HotbitTokenHTB"
    "B21" = "BitKan"
    "C21" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "E21" = "20BitKanKAN"
    "B22" = "NitroEx"
    "C22" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "E22" = "21NitroExNTX"
    "B23" = "LEO"
    "C23" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "E23" = "22LEOLEO"
    "B24" = "BTSEToken"
    "C24" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "E24" = "23BTSETokenBTSE"
    "B25" = "BitpandaEcosystemToken"
    "C25" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "E25" = "24BitpandaEcosystemTokenBEST"
    "B26" = "ProBitToken"
    "C26" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "E26" = "25ProBitTokenPROB"
    "E49" = "48BOLOBOLO"
}

# Numeric-looking "Price" cells (column D) - must stay stored as TEXT (inline string)
# exactly as the source workbook does, so formats like trailing zeros ("3.580",
# "0.3300") are preserved instead of being normalized away by numeric conversion.
$priceUpdates = @{
    "D2" = "248.52"
    "D3" = "21.66"
    "D4" = "5.353"
    "D5" = "0.05611"
    "D6" = "3.419"
    "D7" = "6.398"
    "D8" = "0.8153"
    "D9" = "0.9522"
    "D10" = "0.01155"
    "D11" = "0.1423"
    "D12" = "0.07540"
    "D13" = "0.03211"
    "D14" = "0.03098"
    "D15" = "0.09316"
    "D16" = "3.593"
    "D17" = "0.001592"
    "D18" = "0.04714"
    "D19" = "0.006317"
    "D20" = "0.005069"
    "D21" = "0.001034"
    "D22" = "0.0001501"
    "D23" = "3.775"
    "D24" = "2.147"
    "D25" = "0.3300"
    "D26" = "0.1313"
    "D28" = "0.0003002"
    "D40" = "0.03961"
    "D41" = "0.006975"
    "D42" = "0.1064"
    "D43" = "0.003115"
    "D44" = "0.008778"
    "D45" = "0.00005609"
    "D46" = "0.00000000751"
    "D47" = "0.0005499"
    "D48" = "0.7806"
    "D49" = "0.1736"
    "D50" = "0.00002102"
    "D51" = "0.01011"
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage (NumberFormat "@") before assigning so the digit-string
    # is kept verbatim rather than being parsed into a float.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    # Restore the default/no-op style so no stray formatting is introduced.
    $cell.Style = "Normal"
}
